$wb = $excel.ActiveWorkbook

# --- Insert "is_targeted list" sheet after "analyte_class list" (i.e. before "library_layout list") ---
$afterAnalyte = $wb.Worksheets.Item("analyte_class list")
$isTargetedSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterAnalyte)
$isTargetedSheet.Name = "is_targeted list"
# Build the text values as formulas first (so they are not auto-coerced to booleans),
# then convert them to plain text values in place with a copy / paste-values round-trip.
$isTargetedSheet.Cells.Item(1,1).Formula = "=""TRUE"""
$isTargetedSheet.Cells.Item(2,1).Formula = "=""FALSE"""
$isTargetedRange = $isTargetedSheet.Range("A1:A2")
$isTargetedRange.Copy()
$isTargetedRange.PasteSpecial(-4163)

# --- Insert "is_technical_replicate list" sheet after "library_layout list" (i.e. before "library_final_yield_unit list") ---
$afterLibraryLayout = $wb.Worksheets.Item("library_layout list")
$isTechRepSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterLibraryLayout)
$isTechRepSheet.Name = "is_technical_replicate list"
$isTechRepSheet.Cells.Item(1,1).Formula = "=""TRUE"""
$isTechRepSheet.Cells.Item(2,1).Formula = "=""FALSE"""
$isTechRepRange = $isTechRepSheet.Range("A1:A2")
$isTechRepRange.Copy()
$isTechRepRange.PasteSpecial(-4163)

$excel.CutCopyMode = $false

# --- Update data validation on the main "Export as TSV" sheet ---
$ws = $wb.Worksheets.Item("Export as TSV")

# is_targeted column (N): switch from inline TRUE/FALSE list to reference to new list sheet
$rangeN = $ws.Range("N2:N1048576")
$rangeN.Validation.Modify(3, 1, 1, "='is_targeted list'!`$A`$1:`$A`$2")
$rangeN.Validation.ErrorTitle = "Value must come from list"
$rangeN.Validation.ErrorMessage = "Value must be one of: TRUE / FALSE."

# is_technical_replicate column (V): switch from inline TRUE/FALSE list to reference to new list sheet
$rangeV = $ws.Range("V2:V1048576")
$rangeV.Validation.Modify(3, 1, 1, "='is_technical_replicate list'!`$A`$1:`$A`$2")
$rangeV.Validation.ErrorTitle = "Value must come from list"
$rangeV.Validation.ErrorMessage = "Value must be one of: TRUE / FALSE."
